# Fix up the "Absent" (column H) values so they are consistent with the
# "Total Attendance Count" (column D) values: Absent = 1 when the student
# had no attendance that day (D = 0), otherwise Absent = 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

# Data rows start at row 3 (row 1 = header, row 2 = roll/name info row).
for ($r = 3; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2   # column D: Total Attendance Count
    if ($total -eq 1) {
        $ws.Cells.Item($r, 8).Value = 0
    } else {
        $ws.Cells.Item($r, 8).Value = 1
    }
}
